# "application" -> "app" when preceded by a framework name.
# Five spots are affected:
#   1. "Generated Ruby on Rails application"           -> "...app"          (simple text swap)
#   2. "Engineered API-driven Backbone.js application " -> "...app" + bookmark + " " (gains the _GoBack split)
#   3. "Architected form-driven AngularJS app" + bookmark + " that persists..." -> merged into one run (loses the _GoBack split)
#   4. "Ruby on Rails application that uses jquery-tubular" -> "...app..."  (simple text swap)
#   5. "Ruby on Rails application where you can compete"    -> "...app..." (simple text swap)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: the AngularJS paragraph currently owns the "_GoBack" bookmark that
# splits "...AngularJS app" from " that persists...". The edit removes that
# split (and the bookmark moves to the Backbone.js paragraph instead).
# ---------------------------------------------------------------------------
$angularPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "Architected form-driven AngularJS") {
        $angularPara = $para
        break
    }
}

if ($angularPara -ne $null) {
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }

    # Force the paragraph's runs to coalesce around the former split point.
    # A no-op text assignment is short-circuited by the engine, so round-trip
    # through a temporary character to force a genuine mutation, then remove it.
    $fullText = $angularPara.Range.Text
    $markerIdx = $fullText.IndexOf("AngularJS app")
    if ($markerIdx -ge 0) {
        $mStart = $angularPara.Range.Start + $markerIdx
        $mEnd = $angularPara.Range.Start + $fullText.Length - 1
        $mRange = $d.Range($mStart, $mEnd)
        $original = $mRange.Text
        $mRange.Text = $original + "#"
        $tailRange = $d.Range($mEnd, $mEnd + 1)
        $tailRange.Text = ""
    }
}

# ---------------------------------------------------------------------------
# Step 2: the Backbone.js paragraph gains the split that AngularJS lost:
# "...Backbone.js application " -> "...Backbone.js app" + bookmark + " ".
# ---------------------------------------------------------------------------
$backbonePara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "Engineered API-driven Backbone\.js application") {
        $backbonePara = $para
        break
    }
}

if ($backbonePara -ne $null) {
    $fullText = $backbonePara.Range.Text
    $wordIdx = $fullText.IndexOf("application ")
    if ($wordIdx -ge 0) {
        $appStart = $backbonePara.Range.Start + $wordIdx
        $appEnd = $appStart + "application".Length

        $wordRange = $d.Range($appStart, $appEnd)
        $wordRange.Text = "app"

        $newAppEnd = $appStart + 3
        $bmRange = $d.Range($newAppEnd, $newAppEnd)
        if ($d.Bookmarks.Exists("_GoBack")) {
            $d.Bookmarks("_GoBack").Delete()
        }
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}

# ---------------------------------------------------------------------------
# Step 3: plain "application" -> "app" swaps (no bookmark involved).
# ---------------------------------------------------------------------------
$simpleReplacements = @(
    @{ Old = "Generated Ruby on Rails application using Rails Composer"; New = "Generated Ruby on Rails app using Rails Composer" },
    @{ Old = "Ruby on Rails application that uses jquery-tubular"; New = "Ruby on Rails app that uses jquery-tubular" },
    @{ Old = "Ruby on Rails application where you can compete"; New = "Ruby on Rails app where you can compete" }
)

foreach ($rep in $simpleReplacements) {
    $d.Content.Find.Execute($rep.Old, $true, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2) | Out-Null
}

Write-Output "Done."
